# Updated cryptos list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.372.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.622.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.486"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.68%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.847.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.628.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.365.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0524"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.01%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.18%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.784"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.759.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.409"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  +1.46%  "
